$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.327.73'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.216.01'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.99%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.215.06'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.25%  '
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('E10').Value = '  -1.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.70'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.43%  '
$ws.Range('E12').Value = '  -3.61%  '
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.28'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.744.01'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.428.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.216.93'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.25'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.58%  '
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '507.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.728'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.55'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.99%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.143'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +56.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.97'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.58%  '
$ws.Range('E30').Value = '  -1.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.80%  '
$ws.Range('E32').Value = '  -2.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '28.21'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('E35').Value = '  -5.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.41'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.46'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '500.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0771'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.61%  '
$ws.Range('E40').Value = '  -2.01%  '
$ws.Range('E41').Value = '  +2.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.02'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.71'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.294'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.919.70'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.43'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.03'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.117'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.02%  '
